$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 2422.766666666667, 2679, 2268, 0.03091150124867757),
    @(1, 2425, 2552, 2364, 0.02969153722127279),
    @(2, 2287.866666666667, 2411, 2171, 0.03312266667683919),
    @(3, 2155.466666666667, 2234, 2080, 0.03025951385498047),
    @(4, 2157.3, 2251, 1922, 0.03366438547770182),
    @(5, 1786.433333333333, 1949, 1614, 0.0334805170694987),
    @(6, 2599.766666666667, 2811, 2395, 0.03327964941660563),
    @(7, 2434.8, 2550, 2294, 0.03317662080128988),
    @(8, 2456.833333333333, 2761, 2082, 0.03358786106109619),
    @(9, 2268.2, 2326, 2230, 0.03098969459533691)
)

$row = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $row++
}
